$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 142.84616
$ws.Range("I6").Value = 142.84616
$ws.Range("K6").Value = 428.53848
$ws.Range("M6").Value = -316.53848
$ws.Range("H28").Value = 700.35297
$ws.Range("I28").Value = 600.7857
$ws.Range("K28").Value = 600.7857
$ws.Range("M28").Value = -115.7857
$ws.Range("H43").Value = 11149.8
$ws.Range("J43").Value = 8937.25
$ws.Range("L43").Value = 8937.25
$ws.Range("N43").Value = -9075.25
$ws.Range("H62").Value = 5818.143
$ws.Range("I62").Value = 7336.6665
$ws.Range("K62").Value = 7336.6665
$ws.Range("M62").Value = -6712.6665
$ws.Range("H65").Value = 5818.143
$ws.Range("I65").Value = 7336.6665
$ws.Range("K65").Value = 36683.3325
$ws.Range("M65").Value = -33563.3325
$ws.Range("H69").Value = 5224.75
$ws.Range("I69").Value = 4949.5
$ws.Range("J69").Value = 5500
$ws.Range("K69").Value = 14848.5
$ws.Range("L69").Value = 16500
$ws.Range("M69").Value = -13974.5
$ws.Range("N69").Value = -18248
$ws.Range("H72").Value = 5224.75
$ws.Range("I72").Value = 4949.5
$ws.Range("J72").Value = 5500
$ws.Range("K72").Value = 44545.5
$ws.Range("L72").Value = 49500
$ws.Range("M72").Value = -40177.5
$ws.Range("N72").Value = -58236
$ws.Range("H98").Value = 799.3570999999999
$ws.Range("I98").Value = 822.46155
$ws.Range("K98").Value = 822.46155
$ws.Range("M98").Value = 675.53845
$ws.Range("H100").Value = 1980.1786
$ws.Range("J100").Value = 2766.6667
$ws.Range("L100").Value = 2766.6667
$ws.Range("N100").Value = -3848.6667
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()
$ws.Range("H122").Value = 799.3570999999999
$ws.Range("I122").Value = 822.46155
$ws.Range("K122").Value = 2467.38465
$ws.Range("M122").Value = -17.38464999999997

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1462
$ws.Range("I45").Value = 1243.7142
$ws.Range("K45").Value = 1243.7142
$ws.Range("M45").Value = -866.7141999999999
$ws.Range("H74").Value = 1818.64
$ws.Range("I74").Value = 1816
$ws.Range("K74").Value = 1816
$ws.Range("M74").Value = -942
$ws.Range("H77").Value = 1818.64
$ws.Range("I77").Value = 1816
$ws.Range("K77").Value = 9080
$ws.Range("M77").Value = -4712
$ws.Range("H110").Value = 1571.2727
$ws.Range("I110").Value = 1599.125
$ws.Range("K110").Value = 1599.125
$ws.Range("M110").Value = 445.875
$ws.Range("H125").Value = 59999
$ws.Range("J125").Value = 59999
$ws.Range("L125").Value = 59999
$ws.Range("N125").Value = -69839

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H124").Value = 50000
$ws.Range("J124").Value = 50000
$ws.Range("L124").Value = 50000
$ws.Range("N124").Value = -59820

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 847.4286
$ws.Range("I16").Value = 586.3333
$ws.Range("J16").Value = 1043.25
$ws.Range("K16").Value = 586.3333
$ws.Range("L16").Value = 1043.25
$ws.Range("M16").Value = -299.3333
$ws.Range("N16").Value = -1617.25
$ws.Range("H31").Value = 3024.7273
$ws.Range("I31").Value = 2994
$ws.Range("K31").Value = 2994
$ws.Range("M31").Value = -2699
$ws.Range("H34").Value = 3024.7273
$ws.Range("I34").Value = 2994
$ws.Range("K34").Value = 2994
$ws.Range("M34").Value = -2792
$ws.Range("H74").Value = 33000
$ws.Range("J74").Value = 39500
$ws.Range("L74").Value = 39500
$ws.Range("N74").Value = -41248
$ws.Range("H77").Value = 33000
$ws.Range("J77").Value = 39500
$ws.Range("L77").Value = 118500
$ws.Range("N77").Value = -127236
$ws.Range("H99").Value = 6372.85
$ws.Range("I99").Value = 5715.6
$ws.Range("K99").Value = 5715.6
$ws.Range("M99").Value = -4217.6
$ws.Range("H113").Value = 847.4286
$ws.Range("I113").Value = 586.3333
$ws.Range("J113").Value = 1043.25
$ws.Range("K113").Value = 586.3333
$ws.Range("L113").Value = 1043.25
$ws.Range("M113").Value = 1583.6667
$ws.Range("N113").Value = -5383.25
$ws.Range("H126").Value = 6372.85
$ws.Range("I126").Value = 5715.6
$ws.Range("K126").Value = 17146.8
$ws.Range("M126").Value = -14676.8
$ws.Range("H134").Value = 2087.9167
$ws.Range("I134").Value = 1306.4
$ws.Range("K134").Value = 3919.2
$ws.Range("M134").Value = -1384.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 1000
$ws.Range("I97").Value = 1000
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 3000
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -2504
$ws.Range("N97").ClearContents()
$ws.Range("H119").Value = 2000
$ws.Range("I119").Value = 2000
$ws.Range("K119").Value = 6000
$ws.Range("M119").Value = -1162
$ws.Range("H128").Value = 499992
$ws.Range("I128").Value = 499992
$ws.Range("K128").Value = 1499976
$ws.Range("M128").Value = -1494996
$ws.Range("H131").Value = 1333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 3000
$ws.Range("J23").Value = 3000
$ws.Range("L23").Value = 3000
$ws.Range("N23").Value = -3446
$ws.Range("H107").Value = 863.2308
$ws.Range("I107").Value = 793.4545000000001
$ws.Range("J107").Value = 1247
$ws.Range("K107").Value = 793.4545000000001
$ws.Range("L107").Value = 1247
$ws.Range("M107").Value = 1126.5455
$ws.Range("N107").Value = -5087

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5000
$ws.Range("I7").Value = 5000
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 5000
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -4888
$ws.Range("N7").ClearContents()
$ws.Range("H46").Value = 3250.1
$ws.Range("I46").Value = 1499
$ws.Range("K46").Value = 1499
$ws.Range("M46").Value = -1311
$ws.Range("H93").Value = 2733.3333
$ws.Range("I93").Value = 2600
$ws.Range("J93").Value = 3000
$ws.Range("K93").Value = 2600
$ws.Range("L93").Value = 3000
$ws.Range("M93").Value = -1352
$ws.Range("N93").Value = -5496
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("H126").Value = 5000
$ws.Range("I126").Value = 5000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 15000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -12530
$ws.Range("N126").ClearContents()
$ws.Range("H136").Value = 4989.5
$ws.Range("I136").Value = 4332.5
$ws.Range("J136").Value = 5975
$ws.Range("K136").Value = 12997.5
$ws.Range("L136").Value = 17925
$ws.Range("M136").Value = -10447.5
$ws.Range("N136").Value = -23025

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 839.8
$ws.Range("I81").Value = 799.75
$ws.Range("J81").Value = 1000
$ws.Range("K81").Value = 1599.5
$ws.Range("L81").Value = 2000
$ws.Range("M81").Value = -538.5
$ws.Range("N81").Value = -4122
$ws.Range("H84").Value = 839.8
$ws.Range("I84").Value = 799.75
$ws.Range("J84").Value = 1000
$ws.Range("K84").Value = 7997.5
$ws.Range("L84").Value = 10000
$ws.Range("M84").Value = -2693.5
$ws.Range("N84").Value = -20608
$ws.Range("H132").Value = 4031.375
$ws.Range("I132").Value = 4132.6665
$ws.Range("J132").Value = 3727.5
$ws.Range("K132").Value = 12397.9995
$ws.Range("L132").Value = 11182.5
$ws.Range("M132").Value = -9867.999500000002
$ws.Range("N132").Value = -16242.5
